$wb = $excel.ActiveWorkbook

# --- Settings sheet: insert a new row for the Trello workspace setting ---
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Select() | Out-Null

$wsSettings.Rows.Item(5).Insert() | Out-Null
# Restore the standard row height/format lost on insert so the new row
# matches the sheet's default data-row formatting.
$wsSettings.Rows.Item(5).RowHeight = 14.25

$wsSettings.Cells.Item(5, 1).Value = "TrelloWorkspace"
$wsSettings.Cells.Item(5, 3).Value = "Id of Trello workspace containing cohort boards."
$wsSettings.Cells.Item(5, 2).Value = "userworkspace95961868"

$wsSettings.Range("B15").Select() | Out-Null

# --- Assets sheet: add Trello API key/secret orchestrator assets ---
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Select() | Out-Null

$wsAssets.Cells.Item(2, 1).Value = "TrelloAPIKey"
$wsAssets.Cells.Item(3, 1).Value = "TrelloAPISecret"

$wsAssets.Cells.Item(2, 2).Value = "Trello API Key"
$wsAssets.Cells.Item(3, 2).Value = "Trello API Secret"

$wsAssets.Cells.Item(2, 3).Value = "P3 Automation"
$wsAssets.Cells.Item(3, 3).Value = "P3 Automation"

$wsAssets.Cells.Item(3, 4).Value = "Oauth 1 secret for Trello API."
$wsAssets.Cells.Item(2, 4).Value = "Developer key for Trello API."

$wsAssets.Range("C7").Select() | Out-Null
